# Adds "Total" and "Devided" columns (D & E) to the "4x4" sheet, computing
# the row total (sum of Males + Females) and that total rounded to the
# nearest 10,000 ("devided" by 10000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4x4")
# ($ws is also reachable as $wb.ActiveSheet - the "4x4" tab is the one
# selected in the saved workbook)

# Header row
$ws.Range("D1").Value = "Total"
$ws.Range("E1").Value = "Devided"

# Data rows 2-5: D = SUM(B:C), E = ROUND(D/10000,0)
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=SUM(B${r}:C${r})"
    $ws.Cells.Item($r, 5).Formula = "=ROUND(D${r}/10000,0)"
}

# Move the active selection like in the final workbook
$ws.Range("H15").Select()

$wb.Save()
